$wb = $excel.ActiveWorkbook
$sheetA = $wb.Worksheets.Item(1)   # will become "总计" (A1:D2)
$sheetB = $wb.Worksheets.Item(2)   # will become "2021-Q1" (A1:H12)

# Clear existing contents before re-writing so no stale cells remain.
$sheetA.Cells.Clear()
$sheetB.Cells.Clear()

# Rename through a temporary name to avoid a transient name clash.
$sheetB.Name = "__tmp_swap__"
$sheetA.Name = "总计"
$sheetB.Name = "2021-Q1"

# --- "总计" sheet data (now in $sheetA, positioned first) ---
# B2 holds a quarter label that must stay text, not be parsed as a date/formula.
$sheetA.Range("B2:B2").NumberFormat = "@"

$sheetA.Range("B1").Value = "日期"
$sheetA.Range("C1").Value = "持有数量(只)"
$sheetA.Range("D1").Value = "持有市值(亿元)"
$sheetA.Range("A2").Value = 0
$sheetA.Range("B2").Value = "2021-Q1"
$sheetA.Range("C2").Value = 11
$sheetA.Range("D2").Value = 5.42

# Style: bold, thin box border, centered/top-aligned -- header row + index column A
foreach ($totStyleRange in @($sheetA.Range("B1:D1"), $sheetA.Range("A2:A2"))) {
    $totStyleRange.Font.Bold = $true
    $totStyleRange.HorizontalAlignment = -4108
    $totStyleRange.VerticalAlignment = -4160
    $totStyleRange.Borders.LineStyle = 1
}

# Page margins travel with the "总计" sheet's own layout (0.75/1.0/0.5 in).
$sheetA.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$sheetA.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$sheetA.PageSetup.TopMargin = $excel.InchesToPoints(1)
$sheetA.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$sheetA.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$sheetA.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# --- "2021-Q1" sheet data (now in $sheetB, positioned second) ---
# Columns B:G are text in the source data (fund codes with leading zeros, and
# numeric-looking percentages/amounts stored as text) -- force text format first
# so entering the values does not re-interpret them as numbers.
$sheetB.Range("B2:G12").NumberFormat = "@"

$sheetB.Range("B1").Value = "基金代码"
$sheetB.Range("C1").Value = "基金名称"
$sheetB.Range("D1").Value = "基金金额"
$sheetB.Range("E1").Value = "股票总仓位"
$sheetB.Range("F1").Value = "仓位占比"
$sheetB.Range("G1").Value = "持有市值(亿元)"
$sheetB.Range("H1").Value = "仓位排名"

$sheetB.Range("A2").Value = 0
$sheetB.Range("B2").Value = "270023"
$sheetB.Range("C2").Value = "广发全球精选股票(QDII)"
$sheetB.Range("D2").Value = "50.63"
$sheetB.Range("E2").Value = "89.90"
$sheetB.Range("F2").Value = "3.11"
$sheetB.Range("G2").Value = "1.5746"
$sheetB.Range("H2").Value = 8

$sheetB.Range("A3").Value = 1
$sheetB.Range("B3").Value = "000906"
$sheetB.Range("C3").Value = "广发全球精选股票(QDII)美元现汇"
$sheetB.Range("D3").Value = "50.63"
$sheetB.Range("E3").Value = "89.90"
$sheetB.Range("F3").Value = "3.11"
$sheetB.Range("G3").Value = "1.5746"
$sheetB.Range("H3").Value = 8

$sheetB.Range("A4").Value = 2
$sheetB.Range("B4").Value = "000041"
$sheetB.Range("C4").Value = "华夏全球精选股票(QDII)"
$sheetB.Range("D4").Value = "30.43"
$sheetB.Range("E4").Value = "83.78"
$sheetB.Range("F4").Value = "4.11"
$sheetB.Range("G4").Value = "1.2507"
$sheetB.Range("H4").Value = 4

$sheetB.Range("A5").Value = 3
$sheetB.Range("B5").Value = "001668"
$sheetB.Range("C5").Value = "汇添富全球移动互联灵活配置混合（QDII）"
$sheetB.Range("D5").Value = "27.65"
$sheetB.Range("E5").Value = "92.28"
$sheetB.Range("F5").Value = "2.20"
$sheetB.Range("G5").Value = "0.6083"
$sheetB.Range("H5").Value = 8

$sheetB.Range("A6").Value = 4
$sheetB.Range("B6").Value = "006792"
$sheetB.Range("C6").Value = "鹏华香港美国互联网股票（LOF）美元现汇"
$sheetB.Range("D6").Value = "2.81"
$sheetB.Range("E6").Value = "81.10"
$sheetB.Range("F6").Value = "3.87"
$sheetB.Range("G6").Value = "0.1087"
$sheetB.Range("H6").Value = 4

$sheetB.Range("A7").Value = 5
$sheetB.Range("B7").Value = "160644"
$sheetB.Range("C7").Value = "鹏华香港美国互联网股票（LOF）人民币"
$sheetB.Range("D7").Value = "2.81"
$sheetB.Range("E7").Value = "81.10"
$sheetB.Range("F7").Value = "3.87"
$sheetB.Range("G7").Value = "0.1087"
$sheetB.Range("H7").Value = 4

$sheetB.Range("A8").Value = 6
$sheetB.Range("B8").Value = "378006"
$sheetB.Range("C8").Value = "上投摩根全球新兴市场混合(QDII)"
$sheetB.Range("D8").Value = "0.65"
$sheetB.Range("E8").Value = "90.97"
$sheetB.Range("F8").Value = "8.27"
$sheetB.Range("G8").Value = "0.0538"
$sheetB.Range("H8").Value = 2

$sheetB.Range("A9").Value = 7
$sheetB.Range("B9").Value = "006373"
$sheetB.Range("C9").Value = "富兰克林国海全球科技互联混合（QDII）人民币"
$sheetB.Range("D9").Value = "0.90"
$sheetB.Range("E9").Value = "85.34"
$sheetB.Range("F9").Value = "5.62"
$sheetB.Range("G9").Value = "0.0506"
$sheetB.Range("H9").Value = 1

$sheetB.Range("A10").Value = 8
$sheetB.Range("B10").Value = "006374"
$sheetB.Range("C10").Value = "富兰克林国海全球科技互联混合（QDII）美元现汇"
$sheetB.Range("D10").Value = "0.90"
$sheetB.Range("E10").Value = "85.34"
$sheetB.Range("F10").Value = "5.62"
$sheetB.Range("G10").Value = "0.0506"
$sheetB.Range("H10").Value = 1

$sheetB.Range("A11").Value = 9
$sheetB.Range("B11").Value = "005698"
$sheetB.Range("C11").Value = "华夏全球科技先锋混合QDII"
$sheetB.Range("D11").Value = "0.84"
$sheetB.Range("E11").Value = "82.12"
$sheetB.Range("F11").Value = "3.35"
$sheetB.Range("G11").Value = "0.0281"
$sheetB.Range("H11").Value = 10

$sheetB.Range("A12").Value = 10
$sheetB.Range("B12").Value = "006555"
$sheetB.Range("C12").Value = "浦银安盛全球智能科技股票（QDII）"
$sheetB.Range("D12").Value = "0.52"
$sheetB.Range("E12").Value = "82.10"
$sheetB.Range("F12").Value = "1.87"
$sheetB.Range("G12").Value = "0.0097"
$sheetB.Range("H12").Value = 9

# Style: bold, thin box border, centered/top-aligned -- header row + index column A
foreach ($q1StyleRange in @($sheetB.Range("B1:H1"), $sheetB.Range("A2:A12"))) {
    $q1StyleRange.Font.Bold = $true
    $q1StyleRange.HorizontalAlignment = -4108
    $q1StyleRange.VerticalAlignment = -4160
    $q1StyleRange.Borders.LineStyle = 1
}

# Page margins travel with the "2021-Q1" sheet's own layout (0.7/0.75/0.3 in).
$sheetB.PageSetup.LeftMargin = $excel.InchesToPoints(0.7)
$sheetB.PageSetup.RightMargin = $excel.InchesToPoints(0.7)
$sheetB.PageSetup.TopMargin = $excel.InchesToPoints(0.75)
$sheetB.PageSetup.BottomMargin = $excel.InchesToPoints(0.75)
$sheetB.PageSetup.HeaderMargin = $excel.InchesToPoints(0.3)
$sheetB.PageSetup.FooterMargin = $excel.InchesToPoints(0.3)

# "2021-Q1" keeps the active/selected tab, same as before the resort.
$sheetB.Activate()
$sheetB.Range("A1").Select()